$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, shifting existing rows 132-154 down to 133-155
$ws.Rows(132).Insert()

# Populate the newly inserted row 132 with the new data record
$ws.Range("A132").Value = 4
$ws.Range("B132").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C132").Value = "Los Lagos"
$ws.Range("D132").Value = 44505
$ws.Range("E132").Value = 10
$ws.Range("F132").Value = "Fruta"
$ws.Range("G132").Value = 100108
$ws.Range("H132").Value = "Tropicales y subtropicales"
$ws.Range("I132").Value = 100108005
$ws.Range("J132").Value = "Piña"
$ws.Range("K132").Value = "Caramelo"
$ws.Range("L132").Value = "Segunda"
$ws.Range("M132").Value = 200
$ws.Range("N132").Value = 21000
$ws.Range("O132").Value = 22000
$ws.Range("P132").Value = 21500
$ws.Range("Q132").Value = "$/caja 14 unidades"
$ws.Range("R132").Value = "Ecuador"
$ws.Range("S132").Value = 1536
$ws.Range("T132").Value = 14
